$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 129.14285
$ws.Range("I33").Value = 78.75
$ws.Range("J33").Value = 196.33333
$ws.Range("K33").Value = 78.75
$ws.Range("L33").Value = 196.33333
$ws.Range("M33").Value = 150.25
$ws.Range("N33").Value = -654.3333299999999
$ws.Range("H39").Value = 173.6
$ws.Range("I39").Value = 159.55556
$ws.Range("K39").Value = 478.66668
$ws.Range("M39").Value = -182.66668
$ws.Range("H51").Value = 2604.3333
$ws.Range("I51").Value = 1569
$ws.Range("J51").Value = 2949.4443
$ws.Range("K51").Value = 1569
$ws.Range("L51").Value = 2949.4443
$ws.Range("M51").Value = -1085
$ws.Range("N51").Value = -3917.4443
$ws.Range("H70").Value = 6539.1665
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30540
$ws.Range("H73").Value = 6539.1665
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -31872
$ws.Range("H100").Value = 1004.6667
$ws.Range("I100").Value = 1045.3334
$ws.Range("K100").Value = 1045.3334
$ws.Range("M100").Value = -504.3334
$ws.Range("H125").Value = 606.5
$ws.Range("I125").Value = 642
$ws.Range("K125").Value = 5778
$ws.Range("M125").Value = -3318
$ws.Range("H129").Value = 1998.3334
$ws.Range("I129").Value = 1499.5
$ws.Range("K129").Value = 4498.5
$ws.Range("M129").Value = 501.5
$ws.Range("H132").Value = 12405.375
$ws.Range("I132").Value = 11827.842
$ws.Range("J132").Value = 14600
$ws.Range("K132").Value = 35483.526
$ws.Range("L132").Value = 43800
$ws.Range("M132").Value = -32953.526
$ws.Range("N132").Value = -48860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1481
$ws.Range("I88").Value = 1481
$ws.Range("K88").Value = 1481
$ws.Range("M88").Value = -1075
$ws.Range("H91").Value = 1481
$ws.Range("I91").Value = 1481
$ws.Range("K91").Value = 1481
$ws.Range("M91").Value = -77
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 273.5
$ws.Range("J22").Value = 250
$ws.Range("L22").Value = 250
$ws.Range("N22").Value = -596
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("H107").Value = 5300.1
$ws.Range("I107").Value = 2166.8333
$ws.Range("K107").Value = 2166.8333
$ws.Range("M107").Value = -246.8332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7335.5884
$ws.Range("I31").Value = 2058.2
$ws.Range("K31").Value = 2058.2
$ws.Range("M31").Value = -1763.2
$ws.Range("H34").Value = 7335.5884
$ws.Range("I34").Value = 2058.2
$ws.Range("K34").Value = 2058.2
$ws.Range("M34").Value = -1856.2
$ws.Range("H58").Value = 3308.5833
$ws.Range("I58").Value = 1502.375
$ws.Range("K58").Value = 1502.375
$ws.Range("M58").Value = -1299.375
$ws.Range("H68").Value = 55546.668
$ws.Range("J68").Value = 55546.668
$ws.Range("L68").Value = 55546.668
$ws.Range("N68").Value = -57044.668
$ws.Range("H71").Value = 55546.668
$ws.Range("J71").Value = 55546.668
$ws.Range("L71").Value = 166640.004
$ws.Range("N71").Value = -174128.004
$ws.Range("H99").Value = 2454.2727
$ws.Range("I99").Value = 2454.2727
$ws.Range("K99").Value = 2454.2727
$ws.Range("M99").Value = -956.2727
$ws.Range("H100").Value = 49750
$ws.Range("J100").Value = 49750
$ws.Range("L100").Value = 49750
$ws.Range("N100").Value = -51914
$ws.Range("H119").Value = 30380.5
$ws.Range("J119").Value = 30380.5
$ws.Range("L119").Value = 30380.5
$ws.Range("N119").Value = -40056.5
$ws.Range("H126").Value = 2454.2727
$ws.Range("I126").Value = 2454.2727
$ws.Range("K126").Value = 7362.8181
$ws.Range("M126").Value = -4892.8181
$ws.Range("H132").Value = 3055.8
$ws.Range("I132").Value = 3326
$ws.Range("K132").Value = 9978
$ws.Range("M132").Value = -7448
$ws.Range("H134").Value = 7747.4
$ws.Range("J134").Value = 9949.666999999999
$ws.Range("L134").Value = 29849.001
$ws.Range("N134").Value = -34919.001
$ws.Range("H136").Value = 3308.5833
$ws.Range("I136").Value = 1502.375
$ws.Range("K136").Value = 4507.125
$ws.Range("M136").Value = -1957.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = $null
$ws.Range("H23").Value = 132.66667
$ws.Range("I23").Value = 65.666664
$ws.Range("J23").Value = 199.66667
$ws.Range("K23").Value = 196.999992
$ws.Range("L23").Value = 599.00001
$ws.Range("M23").Value = 38.00000800000001
$ws.Range("N23").Value = -1069.00001
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null
$ws.Range("H55").Value = 5186.615
$ws.Range("I55").Value = 298
$ws.Range("J55").Value = 6075.4546
$ws.Range("K55").Value = 894
$ws.Range("L55").Value = 18226.3638
$ws.Range("M55").Value = -717
$ws.Range("N55").Value = -18580.3638
$ws.Range("H86").Value = 566
$ws.Range("I86").Value = 561.75
$ws.Range("K86").Value = 1685.25
$ws.Range("M86").Value = -499.25
$ws.Range("H89").Value = 566
$ws.Range("I89").Value = 561.75
$ws.Range("K89").Value = 5055.75
$ws.Range("M89").Value = 872.25
$ws.Range("H132").Value = 4566.5
$ws.Range("I132").Value = 4400
$ws.Range("J132").Value = 4599.8
$ws.Range("K132").Value = 39600
$ws.Range("L132").Value = 41398.2
$ws.Range("M132").Value = -37070
$ws.Range("N132").Value = -46458.2
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 811
$ws.Range("I132").Value = 811
$ws.Range("K132").Value = 2433
$ws.Range("M132").Value = 97

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = $null
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = $null
$ws.Range("H43").Value = 28000
$ws.Range("J43").Value = 28000
$ws.Range("L43").Value = 28000
$ws.Range("N43").Value = -28386
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = $null
$ws.Range("H55").Value = 732.5599999999999
$ws.Range("I55").Value = 778.4167
$ws.Range("J55").Value = 690.2308
$ws.Range("K55").Value = 778.4167
$ws.Range("L55").Value = 690.2308
$ws.Range("M55").Value = -605.4167
$ws.Range("N55").Value = -1036.2308
$ws.Range("H82").Value = 3203.6428
$ws.Range("I82").Value = 392.16666
$ws.Range("K82").Value = 392.16666
$ws.Range("M82").Value = -31.16665999999998
$ws.Range("H85").Value = 3203.6428
$ws.Range("I85").Value = 392.16666
$ws.Range("K85").Value = 392.16666
$ws.Range("M85").Value = 855.83334
$ws.Range("H101").Value = 19500
$ws.Range("J101").Value = 19500
$ws.Range("L101").Value = 19500
$ws.Range("N101").Value = -25990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 8001
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null
$ws.Range("H47").Value = 200000
$ws.Range("I47").Value = 200000
$ws.Range("K47").Value = 200000
$ws.Range("M47").Value = -199428
$ws.Range("H100").Value = 1234.7273
$ws.Range("I100").Value = 1234.7273
$ws.Range("K100").Value = 2469.4546
$ws.Range("M100").Value = -1928.4546
$ws.Range("H122").Value = 1313.3334
$ws.Range("I122").Value = 1163.3636
$ws.Range("K122").Value = 3490.0908
$ws.Range("M122").Value = -1040.0908
$ws.Range("H132").Value = 2211.1875
$ws.Range("I132").Value = 2091.9333
$ws.Range("K132").Value = 6275.7999
$ws.Range("M132").Value = -3745.7999
$ws.Range("H136").Value = 3543.6765
$ws.Range("I136").Value = 2899.625
$ws.Range("K136").Value = 8698.875
$ws.Range("M136").Value = -6148.875
